$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 becomes what used to be row 6 (even_MAG-GUT61176.fa)
$ws.Range("A3").Value = "even_MAG-GUT61176.fa"
$ws.Range("B3").Value = 0.02513823199633585
$ws.Range("C3").Value = 0.05767264116206012
$ws.Range("D3").Value = 0.1423179105908852
$ws.Range("E3").Value = 0.26196894358936
$ws.Range("F3").Value = 0.5129022726613588
$ws.Range("G3").Value = 0.5129022726613588
$ws.Range("H3").Value = "s__Limosilactobacillus vaginalis_A"
$ws.Range("I3").Value = "s__Limosilactobacillus vaginalis_A"

# Row 4 becomes what used to be row 10 (even_MAG-GUT83507.fa)
$ws.Range("A4").Value = "even_MAG-GUT83507.fa"
$ws.Range("B4").Value = 0.01605307094748632
$ws.Range("C4").Value = 0.05379831524397638
$ws.Range("D4").Value = 0.1723035000629603
$ws.Range("E4").Value = 0.2934853309602789
$ws.Range("F4").Value = 0.4643597827852982
$ws.Range("G4").Value = 0.4643597827852982
$ws.Range("H4").Value = "s__Limosilactobacillus vaginalis_A"
$ws.Range("I4").Value = "s__Limosilactobacillus vaginalis_A"

# Remove old rows 5-10 (their data has been folded into rows 3 and 4 above)
$ws.Range("A5:I10").Clear()
